$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.977.19"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.384.00"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "3.963.15"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "3.379.79"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "61.065.21"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "3.524.29"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E33").Value = "  -4.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "3.416.42"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0765"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.778"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "2.445.85"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
